# The sector-group lookup table had its "category-name" (column E) and
# "group-code" (column F) columns accidentally transposed. This swaps the
# two columns back for every row (including the header) so that:
#   E = codeforiati:group-code
#   F = codeforiati:category-name

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$lastRowAddr = [string]$lastRow

$colE = $ws.Range("E1:E" + $lastRowAddr)
$colF = $ws.Range("F1:F" + $lastRowAddr)
$scratch = $ws.Range("Z1:Z" + $lastRowAddr)

# Stash column E in an unused scratch column.
$colE.Copy() | Out-Null
$scratch.PasteSpecial() | Out-Null

# Move column F's values into column E.
$colF.Copy() | Out-Null
$colE.PasteSpecial() | Out-Null

# Move the stashed original column E values into column F.
$scratch.Copy() | Out-Null
$colF.PasteSpecial() | Out-Null

# Clean up the scratch column so it doesn't leave any trace in the sheet.
$scratch.ClearContents() | Out-Null
